# Update column F (dSF) values for the affected rows with repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    3  = 3
    5  = 0
    6  = -1
    7  = -1
    8  = 3
    9  = 1
    10 = 1
    11 = 1
    12 = 2
    13 = -5
    14 = 0
    15 = -3
    16 = 9
    17 = 3
    19 = 6
    20 = 2
    21 = 9
    24 = -2
    25 = -1
    26 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
